# MoJ Statistics Forward Look - weekly update 02.10.25
# Applies the changes described by the commit "stats forward look weekly update 02.10.25"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "as at" date in the intro paragraph (A2)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 02 October 2025"

# Week commencing 20 Oct 2025: new publication added (row 6 was previously empty except for the week/date cells)
$ws.Range("B6").Value = "Criminal justice statistics quarterly: June 2025"
$ws.Range("C6").Value = "23 October 2025"
$ws.Range("D6").Value = "confirmed"
$ws.Range("F6").Value = "standard"

# Week commencing 27 Oct 2025: "Justice data lab statistics: October 2025" replaced with a new
# publication, and a new "Offender management statistics quarterly" row inserted before the
# "Deaths of offenders..." row (which also flips to confirmed).
$ws.Range("B7").Value = "Safety in the children and young people secure estate: Update to June 2025"
$ws.Range("B9").Value = "Proven reoffending statistics: October to December 2023"
$ws.Range("B10").Value = "Offender management statistics quarterly: April to June 2025"
$ws.Range("D11").Value = "confirmed"

# Week commencing 24 Nov 2025: "Her Majesty's" -> "His Majesty's", and status flips to confirmed.
$ws.Range("B16").Value = "His Majesty" + [char]8217 + "s Prison and Probation Service offender equalities report: 2024 to 2025"
$ws.Range("D16").Value = "confirmed"
